$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Mid-Term paragraph updates
# ------------------------------------------------------------------
# " exam, scheduled " -> " midterm exam, scheduled "
$d.Content.Find.Execute(" exam, scheduled ", $true, $false, $false, $false, $false, $true, 1, $false, " midterm exam, scheduled ", 2)

# remove " (roughly half way)"
$d.Content.Find.Execute(" (roughly half way)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# "  Details to follow." -> " "
$d.Content.Find.Execute("  Details to follow.", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)

# ------------------------------------------------------------------
# 2. Final/Project paragraph updates
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("The project will include being given a dataset and producing a presentation on your findings", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$r.Text = "There will either be a final project or a midterm"

# Move the (single, document-wide) _GoBack bookmark so that it sits right
# after the text just inserted, matching where the cursor was left after
# the edit. Adding a bookmark with an already-existing name relocates it,
# automatically removing it from its previous location (e.g. after
# "6302: Experimental Statistics II" near the top of the document).
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
